# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1) Metadata!B8 - bump the "Date" value.
# 2) Elements sheet - swap the content of the two mapping columns
#    (AK <-> AL): they were generated in the wrong order ("RIM Mapping"
#    before the new business mapping "Spécification métier vers
#    l'extension ROR AdditionalName"). Swap the header text (row 1) and
#    the per-row values for the rows that actually hold data.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata: Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2) Elements: swap columns AK (37) and AL (38) ---------------------
$ws = $wb.Worksheets.Item("Elements")

# Swap the header row (row 1) text between AK1/AL1.
$hdrAK = $ws.Cells.Item(1, 37).Text
$hdrAL = $ws.Cells.Item(1, 38).Text
$ws.Cells.Item(1, 37).Value = $hdrAL
$ws.Cells.Item(1, 38).Value = $hdrAK

# Swap the data rows that actually differ between AK/AL (rows 2 and 4 are
# blank on both sides, so leave them untouched).
foreach ($r in 3, 5, 6) {
    $valAK = $ws.Cells.Item($r, 37).Text
    $valAL = $ws.Cells.Item($r, 38).Text
    $ws.Cells.Item($r, 37).Value = $valAL
    $ws.Cells.Item($r, 38).Value = $valAK
}
